# "Generate Report for Handoff"
#
# The localization-status report gets a fresh handoff generated for b.md:
# status flips from "Handed back: in sync with en-US" to "Ready for handoff"
# on the Overview sheet, and each locale sheet (zh-cn / de-de) records the
# new handoff xlf file name + handoff datetime for b.md's row.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---- Overview sheet: row for b.md (row 3) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady

# ---- zh-cn sheet: row for b.md (row 3) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusReady
$wsZhCn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-09 00:43:46"

# ---- de-de sheet: row for b.md (row 3) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusReady
$wsDeDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-09 00:43:55"
